$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.628.22"
$ws.Range("E2").Value = "  -3.08%  "

$ws.Range("D3").Value = "3.273.66"
$ws.Range("E3").Value = "  -5.53%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'593.50"
$ws.Range("E5").Value = "  -3.09%  "

$ws.Range("D6").Value = "'151.27"
$ws.Range("E6").Value = "  -9.94%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "3.265.53"
$ws.Range("E8").Value = "  -5.59%  "

$ws.Range("D9").Value = "'0.546"
$ws.Range("E9").Value = "  -8.26%  "

$ws.Range("E10").Value = "  -10.45%  "

$ws.Range("E11").Value = "  -5.09%  "

$ws.Range("E12").Value = "  -10.42%  "

$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "'0.0000247"
$ws.Range("E13").Value = "  -8.48%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'38.57"
$ws.Range("E14").Value = "  -13.35%  "

$ws.Range("D15").Value = "3.795.48"
$ws.Range("E15").Value = "  -5.65%  "

$ws.Range("D16").Value = "67.587.79"

$ws.Range("D17").Value = "3.269.16"
$ws.Range("E17").Value = "  -5.71%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'534.57"
$ws.Range("E18").Value = "  -8.63%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.114"
$ws.Range("E19").Value = "  -5.31%  "

$ws.Range("D20").Value = "'7.15"
$ws.Range("E20").Value = "  -12.74%  "

$ws.Range("D21").Value = "'15.02"
$ws.Range("E21").Value = "  -12.79%  "

$ws.Range("D22").Value = "'0.760"
$ws.Range("E22").Value = "  -10.95%  "

$ws.Range("D23").Value = "'7.91"
$ws.Range("E23").Value = "  -11.66%  "

$ws.Range("D24").Value = "'85.66"
$ws.Range("E24").Value = "  -10.81%  "

$ws.Range("D25").Value = "'13.60"
$ws.Range("E25").Value = "  -10.90%  "

$ws.Range("E27").Value = "  -10.63%  "

$ws.Range("D28").Value = "'8.12"
$ws.Range("E28").Value = "  -6.09%  "

$ws.Range("D29").Value = "'2.17"
$ws.Range("E29").Value = "  -11.77%  "

$ws.Range("D30").Value = "'29.28"
$ws.Range("E30").Value = "  -11.29%  "

$ws.Range("E31").Value = "  -4.58%  "

$ws.Range("E32").Value = "  -5.81%  "

$ws.Range("D33").Value = "'6.66"
$ws.Range("E33").Value = "  -15.65%  "

$ws.Range("D34").Value = "'5.76"
$ws.Range("E34").Value = "  -12.63%  "

$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "'519.95"
$ws.Range("E35").Value = "  -11.51%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("E37").Value = "  -7.37%  "

$ws.Range("D38").Value = "'53.40"
$ws.Range("E38").Value = "  -5.19%  "

$ws.Range("D39").Value = "'0.0860"
$ws.Range("E39").Value = "  -10.64%  "

$ws.Range("D40").Value = "'9.00"
$ws.Range("E40").Value = "  -15.30%  "

$ws.Range("E41").Value = "  -10.02%  "

$ws.Range("D42").Value = "'2.80"
$ws.Range("E42").Value = "  -11.44%  "

$ws.Range("D43").Value = "2.947.11"
$ws.Range("E43").Value = "  -9.39%  "

$ws.Range("D44").Value = "'0.268"
$ws.Range("E44").Value = "  -9.84%  "

$ws.Range("D45").Value = "0.0₃0592"
$ws.Range("E45").Value = "  -15.39%  "

$ws.Range("D46").Value = "'2.20"
$ws.Range("E46").Value = "  -8.76%  "

$ws.Range("D47").Value = "'26.82"
$ws.Range("E47").Value = "  -13.10%  "

$ws.Range("D49").Value = "'2.34"
$ws.Range("E49").Value = "  -16.39%  "

$ws.Range("E50").Value = "  -9.80%  "

$ws.Range("D51").Value = "'123.88"
$ws.Range("E51").Value = "  -7.81%  "
